$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") was regenerated for every data row (2-13): 46062 -> 46063
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# The source data rows got reshuffled; rows 5, 6, 8, 11 and 13 now carry a
# different record (columns A "Beteckning", B "Datum" and G "Area (ha)").
$ws.Range("A5").Value = "A 27724-2022"
$ws.Range("B5").Value = 44743.48386574074
$ws.Range("G5").Value = 1.3

$ws.Range("A6").Value = "A 64445-2023"
$ws.Range("B6").Value = 45280
$ws.Range("G6").Value = 3.7

$ws.Range("A8").Value = "A 46779-2025"
$ws.Range("B8").Value = 45926
$ws.Range("G8").Value = 1.5

$ws.Range("A11").Value = "A 56948-2025"
$ws.Range("B11").Value = 45978.64356481482
$ws.Range("G11").Value = 4.7

$ws.Range("A13").Value = "A 50934-2024"
$ws.Range("B13").Value = 45602
$ws.Range("G13").Value = 0.6
